$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.995.81"
$ws.Range("E2").Value = "  +0.39%  "
$ws.Range("D3").Value = "1.638.28"
$ws.Range("E3").Value = "  -0.22%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  -0.37%  "
$ws.Range("D5").Value = "'214.81"
$ws.Range("E5").Value = "  -0.28%  "
$ws.Range("D6").Value = "'0.5089"
$ws.Range("E6").Value = "  +0.80%  "
$ws.Range("E7").Value = "  -0.34%  "
$ws.Range("D8").Value = "'0.2568"
$ws.Range("E8").Value = "  -0.37%  "
$ws.Range("D9").Value = "'0.06354"
$ws.Range("E9").Value = "  -0.68%  "
$ws.Range("D10").Value = "'19.64"
$ws.Range("E10").Value = "  +0.19%  "
$ws.Range("D11").Value = "'0.07764"
$ws.Range("E11").Value = "  -0.35%  "
$ws.Range("D12").Value = "'4.271"
$ws.Range("E12").Value = "  -0.33%  "
$ws.Range("D13").Value = "1.641.15"
$ws.Range("E13").Value = "  -0.45%  "
$ws.Range("D14").Value = "'0.5437"
$ws.Range("E14").Value = "  +0.05%  "
$ws.Range("D15").Value = "0.0₅7716"
$ws.Range("D16").Value = "'64.01"
$ws.Range("E16").Value = "  -1.41%  "
$ws.Range("D17").Value = "25.999.80"
$ws.Range("E17").Value = "  +0.17%  "
$ws.Range("D18").Value = "'1.002"
$ws.Range("E18").Value = "  -0.38%  "
$ws.Range("D19").Value = "'198.88"
$ws.Range("E19").Value = "  +0.63%  "
$ws.Range("D20").Value = "'4.421"
$ws.Range("E20").Value = "  +0.16%  "
$ws.Range("D21").Value = "'9.910"
$ws.Range("E21").Value = "  -0.54%  "
$ws.Range("D22").Value = "'6.041"
$ws.Range("E22").Value = "  +0.66%  "
$ws.Range("E23").Value = "  -0.39%  "
$ws.Range("D24").Value = "'1.892"
$ws.Range("E24").Value = "  +1.32%  "
$ws.Range("D25").Value = "'141.16"
$ws.Range("E25").Value = "  +0.35%  "
$ws.Range("D26").Value = "'0.1201"
$ws.Range("E26").Value = "  +5.10%  "
$ws.Range("D27").Value = "'6.832"
$ws.Range("E27").Value = "  -0.63%  "
$ws.Range("E28").Value = "  -0.79%  "
$ws.Range("E29").Value = "  -0.91%  "
$ws.Range("D30").Value = "'0.04900"
$ws.Range("E30").Value = "  -2.30%  "
$ws.Range("E31").Value = "  -0.19%  "
$ws.Range("E32").Value = "  -0.90%  "
$ws.Range("D33").Value = "'1.529"
$ws.Range("E33").Value = "  -0.34%  "
$ws.Range("D34").Value = "'2.373"
$ws.Range("E34").Value = "  +0.15%  "
$ws.Range("D35").Value = "'0.9078"
$ws.Range("E35").Value = "  +1.55%  "
$ws.Range("D36").Value = "'2.589"
$ws.Range("E36").Value = "  -0.94%  "
$ws.Range("D37").Value = "1.129.75"
$ws.Range("E37").Value = "  -1.22%  "
$ws.Range("E38").Value = "  -1.66%  "
$ws.Range("D39").Value = "'0.01561"
$ws.Range("E39").Value = "  -0.10%  "
$ws.Range("E40").Value = "  -0.44%  "
$ws.Range("E41").Value = "  -1.68%  "
$ws.Range("D42").Value = "'0.8117"
$ws.Range("E42").Value = "  -1.03%  "
$ws.Range("D43").Value = "0.0₈125"
$ws.Range("E43").Value = "  +1.81%  "
$ws.Range("D44").Value = "'99.06"
$ws.Range("E44").Value = "  -0.94%  "
$ws.Range("D45").Value = "'5.440"
$ws.Range("E45").Value = "  -4.47%  "
$ws.Range("D46").Value = "1.775.52"
$ws.Range("E46").Value = "  -0.33%  "
$ws.Range("D47").Value = "'0.4528"
$ws.Range("E47").Value = "  -0.03%  "
$ws.Range("D48").Value = "'54.93"
$ws.Range("E48").Value = "  -0.85%  "
$ws.Range("D49").Value = "'0.9965"
$ws.Range("E49").Value = "  -1.12%  "
$ws.Range("D50").Value = "'0.05124"
$ws.Range("E50").Value = "  +0.94%  "
$ws.Range("E51").Value = "  -0.24%  "
